$d = $word.ActiveDocument

# Locate the paragraph containing "Quiz 2 (both sections): Thursday, April 1st"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Quiz 2 (both sections): Thursday, April 1st*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Insert a new paragraph immediately after the "Quiz 2" list item. Word
    # automatically carries over the paragraph style and list numbering
    # (pStyle "Compact", numPr ilvl 0 / numId 1002) from the preceding
    # paragraph, matching the other quiz bullet items.
    $target.Range.InsertParagraphAfter()

    $newPara = $target.Next()
    $newRange = $newPara.Range
    # Exclude the trailing paragraph mark when setting the text.
    $textRange = $d.Range($newRange.Start, $newRange.End - 1)
    $textRange.Text = "Quiz 3: to be scheduled by classroom services during the exam period. You must write the exam to pass the course."
}
